# Update country indicator figures on the "Summary" sheet with more
# precise decimal values. The cells hold numeric-looking text (shared
# strings), so we force Text formatting before assigning the new value to
# keep the cell's underlying type as text (matching the original data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B11" = "68.58"
    "C11" = "5.14"
    "D11" = "73.72"
    "B12" = "26.99"
    "C12" = "34.06"
    "D12" = "61.04"
    "B33" = "31.53"
    "C33" = "4.47"
    "B34" = "24.42"
    "C34" = "43.27"
    "D34" = "67.69"
    "B36" = "87.29"
    "C36" = "12.37"
    "D36" = "99.66"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
